$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -4065324268110.38
$ws.Range("C2").Value = 95957396481.2445
$ws.Range("B3").Value = -3989262069558.1
$ws.Range("C3").Value = 94097441744.1495
$ws.Range("B4").Value = -3953306880030.46
$ws.Range("C4").Value = 93151127400.904
$ws.Range("B5").Value = -3949435416956.16
$ws.Range("C5").Value = 92975993985.7815
$ws.Range("B6").Value = -3976886658617.73
$ws.Range("C6").Value = 93562409623.309
$ws.Range("B7").Value = -4009047137892.59
$ws.Range("C7").Value = 94237393729.133
$ws.Range("B8").Value = -4027725483173.1
$ws.Range("C8").Value = 94587726297.549
$ws.Range("B9").Value = -4041168729553.24
$ws.Range("C9").Value = 95043808275.5815
$ws.Range("B10").Value = -4101715479749.33
$ws.Range("C10").Value = 96406038250.177
$ws.Range("B11").Value = -4242795208256.3
$ws.Range("C11").Value = 99674081556.087
$ws.Range("B12").Value = -4326836255302.77
$ws.Range("C12").Value = 101591013985.918
$ws.Range("B13").Value = -4423878969091.55
$ws.Range("C13").Value = 103746233130.458
$ws.Range("B14").Value = -4505245623726.3
$ws.Range("C14").Value = 105638473309.39
$ws.Range("B15").Value = -4573181702437.52
$ws.Range("C15").Value = 107173725483.292
$ws.Range("B16").Value = -4705679055129.85
$ws.Range("C16").Value = 110247457654.318
$ws.Range("B17").Value = -4800147685396.76
$ws.Range("C17").Value = 112420527810.519
$ws.Range("B18").Value = -4814636080939.78
$ws.Range("C18").Value = 112677142709.303
$ws.Range("B19").Value = -4951513423906.7
$ws.Range("C19").Value = 115815197198.74
$ws.Range("B20").Value = -5021284805915.08
$ws.Range("C20").Value = 117356477824.757
$ws.Range("B21").Value = -5219764134727.81
$ws.Range("C21").Value = 121974915677.837
$ws.Range("B22").Value = -5269430999031.91
$ws.Range("C22").Value = 123174844077.627
$ws.Range("B23").Value = -5488971021948.97
$ws.Range("C23").Value = 128289848587.982
$ws.Range("B24").Value = -5668431147660.3
$ws.Range("C24").Value = 132470967676.927
$ws.Range("B25").Value = -5825217982847.4
$ws.Range("C25").Value = 136154774147.781
$ws.Range("B26").Value = -5990318524835.1
$ws.Range("C26").Value = 139977439813.867
$ws.Range("B27").Value = -6113660092960.73
$ws.Range("C27").Value = 142801592428.955
$ws.Range("B28").Value = -6249729872007.34
$ws.Range("C28").Value = 145785965230.404
$ws.Range("B29").Value = -6425018571908.95
$ws.Range("C29").Value = 149988467608.935
$ws.Range("B30").Value = -6514448419323.52
$ws.Range("C30").Value = 152069022199.053
